$d = $word.ActiveDocument

function Replace-RangeWithParaXml($range, $innerXml) {
    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $innerXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $range.InsertXML($pkg)
}

# 1) Rename the bookmark wrapping the "Online resources for information security"
#    Heading2 paragraph from "online-resources-for-information-security" to
#    "Xeba3c2f3420c339c3ff5bf9459ba3741edb2a9a". The Bookmarks collection's
#    Name setter/Delete are unreliable in this host, so rebuild the paragraph's
#    XML in place (using the Document's own Paragraphs collection, which keeps
#    the trailing paragraph mark in the range so InsertXML performs a clean
#    replace instead of an insert).
$headingParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "Online resources for information security") {
        $headingParaIndex = $i
        break
    }
}

if ($headingParaIndex -ne -1) {
    $headingRange = $d.Paragraphs($headingParaIndex).Range
    $headingXml = '<w:p><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:bookmarkStart w:id="25" w:name="Xeba3c2f3420c339c3ff5bf9459ba3741edb2a9a"/><w:r><w:t xml:space="preserve">Online resources for information security</w:t></w:r><w:bookmarkEnd w:id="25"/></w:p>'
    Replace-RangeWithParaXml $headingRange $headingXml
}

# 2) In the "Class meetings" table, the "Due" column cells that have no text
#    currently contain an empty paragraph styled "Compact"; strip the style so
#    they become plain empty paragraphs.
$tbl = $d.Tables.Item(1)
for ($r = 1; $r -le $tbl.Rows.Count; $r++) {
    $cell = $tbl.Cell($r, 3)
    if ($cell.Range.Text.Trim([char]13, [char]7) -eq "") {
        $p = $cell.Range.Paragraphs.Item(1)
        Replace-RangeWithParaXml $p.Range '<w:p/>'
    }
}

Write-Host "Edits applied"
